$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value  = "https://www.lanacion.com.ar/deportes/futbol/river-volvio-a-ganarle-a-boca-en-reserva-y-disfruta-de-un-ano-perfecto-en-el-superclasico-nid08092024/"
$ws.Range("A8").Value  = "https://www.lanacion.com.ar/deportes/tenis/jannik-sinner-campeon-del-us-open-el-italiano-vencio-a-taylor-fritz-y-se-encumbro-tras-el-ruido-que-nid08092024/"
$ws.Range("A9").Value  = "https://www.lanacion.com.ar/el-mundo/desde-espana-edmundo-gonzalez-urrutia-denuncio-presiones-coacciones-y-amenazas-para-impedir-su-nid08092024/"
$ws.Range("A10").Value = "https://www.lanacion.com.ar/seguridad/asesinaron-a-una-sargento-de-la-policia-bonaerense-mientras-trabajaba-como-chofer-de-una-app-nid08092024/"
$ws.Range("A11").Value = "https://www.lanacion.com.ar/politica/bandera-roja-y-martillo-la-justicia-ordeno-retomar-los-remates-de-bienes-de-lazaro-baez-nid08092024/"

$ws.Range("A16").Value = "https://www.lanacion.com.ar/estados-unidos/donald-trump-amenaza-con-encarcelar-a-quienes-cometieron-fraude-electoral-en-2020-nid08092024/"

$ws.Range("A18").Value = "https://www.lanacion.com.ar/el-mundo/un-patrullero-atropello-a-tres-personas-durante-una-persecucion-en-un-exclusivo-barrio-de-madrid-nid08092024/"

$ws.Range("A20").Value = "https://www.lanacion.com.ar/politica/los-diplomaticos-van-al-paro-y-apuntan-a-mondino-por-el-pago-de-ganancias-nid08092024/"
$ws.Range("A21").Value = "https://www.lanacion.com.ar/politica/bandera-roja-y-martillo-la-justicia-ordeno-retomar-los-remates-de-bienes-de-lazaro-baez-nid08092024/"

$ws.Range("A24").Value = "https://www.lanacion.com.ar/el-mundo/colombia-o-argentina-inteligencia-artificial-predice-resultado-del-esperado-partido-nid08092024/"
$ws.Range("A25").Value = "https://www.lanacion.com.ar/sociedad/me-enseno-a-posicionar-mi-negocio-el-fuerte-impacto-de-un-programa-educativo-gratuito-y-con-nid08092024/"
$ws.Range("A26").Value = "https://www.lanacion.com.ar/economia/martinez-sosa-cobro-comisiones-de-200-mil-polizas-optativas-con-una-trampa-en-gendarmeria-nid08092024/"
$ws.Range("A27").Value = "https://www.lanacion.com.ar/espectaculos/personajes/pacho-odonnell-el-edadismo-su-adiccion-a-la-gimnasia-su-regreso-a-los-escenarios-y-por-que-hay-que-nid08092024/"

$ws.Range("A29").Value = "https://www.lanacion.com.ar/deportes/atletismo/el-increible-accidente-de-la-paratleta-italiana-que-perdio-la-medalla-dorada-y-conmovio-a-todos-en-nid08092024/"
$ws.Range("A30").Value = "https://www.lanacion.com.ar/deportes/futbol/el-golazo-de-chilena-de-dinamarca-en-la-victoria-ante-serbia-por-la-nations-league-la-pirueta-de-nid08092024/"
$ws.Range("A31").Value = "https://www.lanacion.com.ar/el-mundo/maria-corina-machado-resiste-el-hostigamiento-edmundo-luchara-desde-afuera-y-volvera-para-asumir-la-nid08092024/"

$ws.Range("A33").Value = "https://www.lanacion.com.ar/el-mundo/la-larga-lista-de-opositores-perseguidos-que-sufren-el-exilio-impuesto-por-hugo-chavez-y-nicolas-nid08092024/"

$ws.Range("A35").Value = "https://www.lanacion.com.ar/sociedad/revelador-estudio-uno-de-cada-cuatro-pacientes-en-coma-sigue-teniendo-algun-tipo-de-consciencia-nid08092024/"
$ws.Range("A36").Value = "https://www.lanacion.com.ar/el-mundo/londres-38-la-historia-que-cruza-a-pinochet-con-un-jerarca-nazi-fugado-a-chile-nid08092024/"
$ws.Range("A37").Value = "https://www.lanacion.com.ar/seguridad/adn-del-crimen-confirman-la-condena-para-el-poderoso-juez-que-era-la-pieza-clave-de-la-banda-formada-nid08092024/"
$ws.Range("A38").Value = "https://www.lanacion.com.ar/espectaculos/musica/conociendo-rusia-como-salvar-el-rock-con-buenas-canciones-una-guitarra-y-todo-el-carisma-nid08092024/"
$ws.Range("A39").Value = "https://www.lanacion.com.ar/deportes/futbol/la-jugada-del-arquero-de-penarol-que-desato-un-escandalo-tres-expulsados-y-seis-amonestados-nid08092024/"

$ws.Range("A42").Value = "https://www.lanacion.com.ar/seguridad/un-policia-murio-tras-recibir-una-patada-en-la-espalda-cuando-vecinos-defendieron-a-un-delincuente-nid07092024/"
$ws.Range("A43").Value = "https://www.lanacion.com.ar/espectaculos/musica/waldo-de-los-rios-el-genio-musical-que-cumpliria-90-anos-sus-frustraciones-un-auto-incendiado-y-el-nid07092024/"

$ws.Range("A70").Value = "https://www.lanacion.com.ar/politica/el-pacto-no-escrito-entre-javier-milei-y-cristina-kirchner-nid07092024/"
$ws.Range("A71").Value = "https://www.lanacion.com.ar/opinion/milei-mas-rapido-que-kirchner-para-atacar-a-la-prensa-nid08092024/"
$ws.Range("A72").Value = "https://www.lanacion.com.ar/opinion/un-libertario-esclavo-de-sus-promesas-nid08092024/"

$ws.Range("A75").Value = "https://www.lanacion.com.ar/opinion/la-logica-del-reality-invadio-a-la-politica-nid07092024/"

$ws.Range("A83").Value = "https://www.lanacion.com.ar/sociedad/natividad-de-la-santisima-virgen-maria-su-historia-y-que-oracion-rezar-para-pedir-su-ayuda-nid08092024/"
$ws.Range("A84").Value = "https://www.lanacion.com.ar/politica/internaron-al-jefe-de-gabinete-guillermo-francos-por-un-cuadro-gastrointestinal-agudo-nid08092024/"
$ws.Range("A85").Value = "https://www.lanacion.com.ar/economia/jubilados-anses-fijo-el-haber-maximo-a-partir-de-septiembre-en-cuanto-queda-nid08092024/"
$ws.Range("A86").Value = "https://www.lanacion.com.ar/sociedad/que-es-la-lipotimia-el-cuadro-que-sufrio-guillermo-francos-antes-de-ser-internado-nid08092024/"
